$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "estimation" (B) and "real" (C) time values for the tasks that
# reverted to not-yet-estimated state. The "Task velocity" column (D) is a
# calculated table column (IFERROR formula) so it will recompute to "" on
# its own once B/C are blank.
$rowsToClear = @(6, 7, 11, 12, 13, 14, 15, 22)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).ClearContents()
}

# Remove the last two rows of the table (Tableau1), which drops the
# "upvote des commentaires" and "menu arborescent" tasks entirely.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()

# Restore the previously selected cell.
$ws.Range("I17").Select()

$wb.Save()
